$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new Q&A row
$ws.Range("A2").Value = "what is her name"
$ws.Range("B2").Value = "Riddhi Hedaoo"

# Update the selected cell to match the target state
$ws.Range("B9").Select()
